# "Rewrote the debug() function for Excel service"
#
# Net effect on the workbook (docs/services/spreadsheets/needs2.xlsx):
#   1. The assignee name "Celestina Duodu" is renamed to "Haiyang Zhang"
#      everywhere it appears (column F / "assigned to"): rows 3, 7, 11,
#      15 and 19 all share that one string.
#   2. Rows 7, 11, 15 and 19 additionally pick up a distinct (but
#      visually identical) cell style in column F - a second style entry
#      that is a duplicate of the plain "General" style already used by
#      the rest of the row (xfId 3) - while row 3's F cell keeps its
#      original (unstyled) formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the shared "assigned to" value everywhere it occurs ----
# xlWhole (2) so we only match whole-cell contents, not substrings.
$ws.Cells.Replace("Celestina Duodu", "Haiyang Zhang", 2) | Out-Null

# --- 2. Give F7 / F11 / F15 / F19 their own (duplicate) cell style -----
# Writing an alignment property at its already-default value forces the
# engine to allocate a new xf record for these cells instead of
# continuing to share the xf used by the rest of the row/table, which
# mirrors the extra <xf .../> entry added to cellXfs in the original
# edit, without altering how the cells actually look.
foreach ($addr in @("F7", "F11", "F15", "F19")) {
    $ws.Range($addr).ShrinkToFit = $false
}
